# ETRM MT & Date
# Updates the Date / DateandTime / Time / SequenceNo (and the extra
# trailing column on Sheet1) test-data values across the four sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet1 ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("K2").Value = "30-05-2024"
$ws1.Range("N2").Value = "30-05-2024 02:53:30 PM"
$ws1.Range("O2").Value = "02:35:55 PM"
$ws1.Range("AG2").Value = "ET467"
$ws1.Range("AK2").Value = "3"

# --- Sheet2 ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("K2").Value = "30-05-2024"
$ws2.Range("N2").Value = "30-05-2024 02:59:31 PM"
$ws2.Range("O2").Value = "02:35:55 PM"
$ws2.Range("AG2").Value = "ET468"

# --- Sheet3 ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("K2").Value = "30-05-2024"
$ws3.Range("N2").Value = "30-05-2024 02:59:31 PM"
$ws3.Range("O2").Value = "02:35:55 PM"
$ws3.Range("AG2").Value = "ET469"

# --- Sheet4 ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("K2").Value = "30-05-2024"
$ws4.Range("N2").Value = "30-05-2024 02:59:31 PM"
$ws4.Range("O2").Value = "02:35:55 PM"
$ws4.Range("AG2").Value = "ET469"
